$d = $word.ActiveDocument

$payload = @'
<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Las acciones de Lectura, Actualización, Eliminación y Adición a la base de datos se realizarán mediante la </w:t></w:r><w:r><w:t xml:space="preserve">API </w:t></w:r><w:r><w:t>que se creará para el proyecto, la cual brindará la posibilidad de que sea independiente de la Aplicación Web inicial y se pueda utilizar en futuras implementaciones sin que sea necesaria modificación alguna en lo realizado con respecto al “</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>BackEnd</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>”.</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">La </w:t></w:r><w:r><w:t>API</w:t></w:r><w:r><w:t xml:space="preserve"> contendrá la definición de CRUD (CREATE, READE, UPDATE, DELETE) además de otras funciones necesarias para el correcto funcionamiento de la Aplicación Web o cualquier otro cliente que haga uso del mismo.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$count = $d.Paragraphs.Count
$lastPara = $d.Paragraphs.Item($count)
$lastPara.Range.InsertParagraphAfter()

$newCount = $d.Paragraphs.Count
$newPara = $d.Paragraphs.Item($newCount)
$newPara.Range.InsertXML($payload)

$afterCount = $d.Paragraphs.Count
$tailPara = $d.Paragraphs.Item($afterCount)
$delRange = $d.Range($tailPara.Range.Start - 1, $tailPara.Range.End)
$delRange.Delete()

Write-Output ("ParagraphCount: " + $d.Paragraphs.Count)
